$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 12, shifting the existing data (old rows 12-33)
# down to rows 14-35.
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()

# Common column values shared by every data row in this table.
$colA = 10
$colB = "Vega Modelo de Temuco"
$colC = "La Araucanía"
$colE = 9
$colF = "Fruta"
$colG = 100104
$colH = "Frutos de pepita"
$colI = 100104004
$colJ = "Níspero"
$colK = "Californiana(o)"

# New row 12
$ws.Range("A12").Value = $colA
$ws.Range("B12").Value = $colB
$ws.Range("C12").Value = $colC
$ws.Range("D12").Value = 45246
$ws.Range("E12").Value = $colE
$ws.Range("F12").Value = $colF
$ws.Range("G12").Value = $colG
$ws.Range("H12").Value = $colH
$ws.Range("I12").Value = $colI
$ws.Range("J12").Value = $colJ
$ws.Range("K12").Value = $colK
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 120
$ws.Range("N12").Value = 30000
$ws.Range("O12").Value = 30000
$ws.Range("P12").Value = 30000
$ws.Range("Q12").Value = "$/bandeja 10 kilos"
$ws.Range("R12").Value = "Provincia de Quillota"
$ws.Range("S12").Value = 3000
$ws.Range("T12").Value = 10

# New row 13
$ws.Range("A13").Value = $colA
$ws.Range("B13").Value = $colB
$ws.Range("C13").Value = $colC
$ws.Range("D13").Value = 45246
$ws.Range("E13").Value = $colE
$ws.Range("F13").Value = $colF
$ws.Range("G13").Value = $colG
$ws.Range("H13").Value = $colH
$ws.Range("I13").Value = $colI
$ws.Range("J13").Value = $colJ
$ws.Range("K13").Value = $colK
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 22000
$ws.Range("O13").Value = 22000
$ws.Range("P13").Value = 22000
$ws.Range("Q13").Value = "$/bandeja 5 kilos"
$ws.Range("R13").Value = "Provincia de Quillota"
$ws.Range("S13").Value = 4400
$ws.Range("T13").Value = 5
